$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header column F1 with value "RXNO_DEF", matching the
# existing header cells' formatting (bold, centered, thin-bordered -
# same style as B1:E1) by copying the format from E1, then overwriting
# the value with the new header text.
$f1 = $ws.Range("F1")
$ws.Range("E1").Copy($f1)
$f1.Value = "RXNO_DEF"
